$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The debouncing capacitors C18,C25,C33,C34 used to share the same
# 100nF/0603 BOM line as C10,C12,C14,C15,C29,C30,C31,C32. Split them out
# into their own BOM line with a reduced 680pF value, by inserting a new
# row right after the "1u / 16V / 0603" line (row 5) and before the
# "R1,R5 / 24k" line (old row 6).
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).RowHeight = 15

# New row: the four debouncing caps now use a smaller-value part.
$ws.Range("A6").Value = "C18,C25,C33,C34"
$ws.Range("B6").Value = "0603B681K500NT"
$ws.Range("C6").Value = 4

# Trim the designator list / quantity of the original 100nF line (now
# row 4) since C18,C25,C33,C34 moved out of it.
$ws.Range("A4").Value = "C10,C12,C14,C15,C29,C30,C31,C32"
$ws.Range("C4").Value = 8

# Fill in the new row's value last (matches authoring order).
$ws.Range("D6").Value = "680p / 50V / 0603"

# A handful of other Value cells get a " / 0603" package suffix appended
# for consistency with the rest of the table.
$ws.Range("D18").Value = "2k2 / 0603"
$ws.Range("D19").Value = "2k / 0603"
$ws.Range("D20").Value = "330 / 0603"
$ws.Range("D21").Value = "22k / 0603"

# Restore the cursor / selection position recorded in the saved file.
$ws.Range("G19").Select()
